# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off again (new xliff generated), on the Overview sheet and on the
# per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# NOTE: the "False"/"True" values need a leading apostrophe text-qualifier so
# Excel stores them as plain text (shared strings) instead of native Booleans
# -- this matches the original workbook, where that column already held the
# text strings "True"/"False" rather than boolean cells. The apostrophe
# itself is only a text-qualifier and is not part of the stored cell value.
# Other values (dates, file names, messages) are plain text already and are
# set directly so the cells keep their original number-format style.

# --- Overview sheet: row for b.md (row 3) ---
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-29 22:39:18"

# --- zh-cn sheet: row for b.md (row 3) ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-29 22:39:14"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1afac61636e712d41f546ba5c14c744359fbd546/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc7d4e6f606c313ab923366a887ff6f7fb0a6f08/e2e/b.md."

# --- de-de sheet: row for b.md (row 3) ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-29 22:39:18"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1afac61636e712d41f546ba5c14c744359fbd546/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc7d4e6f606c313ab923366a887ff6f7fb0a6f08/e2e/b.md."

# Widen the "Error Detail" column (P) on both locale sheets so the new long
# message is readable. The engine's ColumnWidth setter adds 5/6 of a
# character to the stored OOXML column width, so subtract that back off to
# land on an OOXML width of exactly 40.
$errorDetailWidth = 40 - 5/6
$zhcn.Columns.Item(16).ColumnWidth = $errorDetailWidth
$dede.Columns.Item(16).ColumnWidth = $errorDetailWidth

Write-Host "Report updated for handoff."
